# "Finita tabella risultati test di accettazione"
# Fill in the specific acceptance-test identifiers ("TA: ...") for each
# user-story card, fix a typo in the "Aggiungere spesa al viaggio" card
# title, wrap the text of the now-longer TA cell for the "Modifica
# viaggio" card, and leave the selection on that cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Registrazione
$ws.Range("A3").Value = "TA: register"
# Visualizza inviti (login)
$ws.Range("G3").Value = "TA: login"

# Visualizza lista viaggi
$ws.Range("A15").Value = "TA: mytravel"
# Modifica viaggio
$ws.Range("G15").Value = "TA: form_modify_travel"

# Manda inviti
$ws.Range("A27").Value = "TA: send_invite"
# Creazione nuovo viaggio
$ws.Range("G27").Value = "TA: create_travel"

# Logout
$ws.Range("A39").Value = "TA: logout"
# Utilizza chat
$ws.Range("G39").Value = "TA: write_comment"

# Visualizza dettagli viaggio
$ws.Range("A51").Value = "TA: details_travel"
# Aggiungere spesa al viaggio
$ws.Range("G51").Value = "TA: form_expense"

# Fix typo in card title: "Aggiugnere" -> "Aggiungere"
$ws.Range("G49").Value = "Titolo: Aggiungere spesa la viaggio"

# The new TA label for "Modifica viaggio" is longer, so enable wrap text
# on its merged cell (G15:H16) like the other long-text card bodies.
$ws.Range("G15:H16").WrapText = $true

# Leave selection on the cell that was just edited.
[void]$ws.Range("G15:H16").Select()
